# Update NATMI TPM LR-pair results for Il1rapl1-Ptprs with the recomputed
# (new TPM-based) numbers: the "ECs -> *" sending-cluster rows are gone and
# the three remaining "MuSCs -> *" rows get refreshed metric columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the three old "ECs" sending-cluster rows (rows 2-4); the three
# "MuSCs" sending-cluster rows shift up to become rows 2-4.
$ws.Range("A2:A4").EntireRow.Delete()

# New values for rows 2-4, columns B..T (column A is already "MuSCs" for
# all three rows after the deletion above, so it is left untouched):
#   B..D = ligand symbol / receptor symbol / target cluster labels
#   E..T = the recomputed metric columns
$rowValues = @{
    2 = @("Il1rapl1", "Ptprs", "ECs",   3, 1, 0.257516, 0.772548, 1, 1, 3, 1, 1.660421,          4.981262999999999, 0.03714789785507311, 0.03714789785507311, 0.4275849742359999, 3.848264768123999, 0.03714789785507311, 0.03714789785507311)
    3 = @("Il1rapl1", "Ptprs", "FAPs",  3, 1, 0.257516, 0.772548, 1, 1, 3, 1, 25.17096033333333, 75.51288099999999, 0.5631392661118858,  0.5631392661118859,  6.481925021198666,  58.33732519078799, 0.5631392661118858,  0.5631392661118859)
    4 = @("Il1rapl1", "Ptprs", "MuSCs", 3, 1, 0.257516, 0.772548, 1, 1, 3, 1, 17.866195,          53.598585,          0.399712836033041,   0.399712836033041,   4.600831071620001,  41.40747964458,    0.399712836033041,   0.399712836033041)
}

foreach ($r in 2..4) {
    $values = $rowValues[$r]
    for ($i = 0; $i -lt $values.Length; $i++) {
        $col = $i + 2   # column B = 2 .. column T = 20
        $ws.Cells.Item($r, $col).Value = $values[$i]
    }
}
